$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 9.65999984741211
$ws.Range("E2").Value = 9.720000267028809
$ws.Range("F2").Value = 9.989999771118164
$ws.Range("G2").Value = 9.65999984741211
$ws.Range("H2").Value = 26730988
$ws.Range("I2").Value = "MDAI"
$ws.Range("D3").Value = 9.65999984741211
$ws.Range("E3").Value = 9.720000267028809
$ws.Range("F3").Value = 9.989999771118164
$ws.Range("G3").Value = 9.65999984741211
$ws.Range("H3").Value = 26730988
$ws.Range("I3").Value = "MDAI"
$ws.Range("D4").Value = 9.68000030517578
$ws.Range("E4").Value = 9.65999984741211
$ws.Range("F4").Value = 9.850000381469728
$ws.Range("G4").Value = 9.630000114440918
$ws.Range("H4").Value = 26730988
$ws.Range("I4").Value = "MDAI"
$ws.Range("D5").Value = 9.710000038146973
$ws.Range("E5").Value = 9.720000267028809
$ws.Range("F5").Value = 9.789999961853027
$ws.Range("G5").Value = 9.670000076293944
$ws.Range("H5").Value = 26730988
$ws.Range("I5").Value = "MDAI"
$ws.Range("D6").Value = 9.8149995803833
$ws.Range("E6").Value = 9.84000015258789
$ws.Range("F6").Value = 9.850000381469728
$ws.Range("G6").Value = 9.789999961853027
$ws.Range("H6").Value = 26730988
$ws.Range("I6").Value = "MDAI"
$ws.Range("D7").Value = 9.76200008392334
$ws.Range("E7").Value = 9.800000190734863
$ws.Range("F7").Value = 9.84000015258789
$ws.Range("G7").Value = 9.760000228881836
$ws.Range("H7").Value = 26730988
$ws.Range("I7").Value = "MDAI"
$ws.Range("D8").Value = 9.8100004196167
$ws.Range("E8").Value = 9.829999923706056
$ws.Range("F8").Value = 9.850000381469728
$ws.Range("G8").Value = 9.800000190734863
$ws.Range("H8").Value = 26730988
$ws.Range("I8").Value = "MDAI"
$ws.Range("D9").Value = 9.850000381469728
$ws.Range("E9").Value = 9.920000076293944
$ws.Range("F9").Value = 9.970000267028809
$ws.Range("G9").Value = 9.850000381469728
$ws.Range("H9").Value = 26730988
$ws.Range("I9").Value = "MDAI"
$ws.Range("D10").Value = 9.949999809265137
$ws.Range("E10").Value = 10.14999961853027
$ws.Range("F10").Value = 10.64999961853027
$ws.Range("G10").Value = 9.949999809265137
$ws.Range("H10").Value = 26730988
$ws.Range("I10").Value = "MDAI"
$ws.Range("D11").Value = 9.970000267028809
$ws.Range("E11").Value = 10.5
$ws.Range("F11").Value = 16.25
$ws.Range("G11").Value = 9.960000038146973
$ws.Range("H11").Value = 26730988
$ws.Range("I11").Value = "MDAI"
$ws.Range("D12").Value = 10.47999954223633
$ws.Range("E12").Value = 10.64000034332275
$ws.Range("F12").Value = 17
$ws.Range("G12").Value = 10.43000030517578
$ws.Range("H12").Value = 26730988
$ws.Range("I12").Value = "MDAI"
$ws.Range("D13").Value = 2.670000076293945
$ws.Range("E13").Value = 2.5
$ws.Range("F13").Value = 4.679999828338623
$ws.Range("G13").Value = 2.200000047683716
$ws.Range("H13").Value = 26730988
$ws.Range("I13").Value = "MDAI"
$ws.Range("D14").Value = 2.5
$ws.Range("E14").Value = 2.089999914169312
$ws.Range("F14").Value = 3.720999956130981
$ws.Range("G14").Value = 1.899999976158142
$ws.Range("H14").Value = 26730988
$ws.Range("I14").Value = "MDAI"
$ws.Range("D15").Value = 2.150000095367432
$ws.Range("E15").Value = 1.779999971389771
$ws.Range("F15").Value = 2.619999885559082
$ws.Range("G15").Value = 1.730000019073486
$ws.Range("H15").Value = 26730988
$ws.Range("I15").Value = "MDAI"
$ws.Range("D16").Value = 1.799999952316284
$ws.Range("E16").Value = 1.669999957084656
$ws.Range("F16").Value = 2.099999904632568
$ws.Range("G16").Value = 1.628999948501587
$ws.Range("H16").Value = 26730988
$ws.Range("I16").Value = "MDAI"
$ws.Range("D17").Value = 1.039999961853027
$ws.Range("E17").Value = 1.440000057220459
$ws.Range("F17").Value = 1.590000033378601
$ws.Range("G17").Value = 0.9610000252723694
$ws.Range("H17").Value = 26730988
$ws.Range("I17").Value = "MDAI"
$ws.Range("D18").Value = 2.880000114440918
$ws.Range("E18").Value = 1.929999947547913
$ws.Range("F18").Value = 3.029999971389771
$ws.Range("G18").Value = 1.720000028610229
$ws.Range("H18").Value = 26730988
$ws.Range("I18").Value = "MDAI"
$ws.Range("D19").Value = 1.169999957084656
$ws.Range("E19").Value = 1.299999952316284
$ws.Range("F19").Value = 1.429999947547913
$ws.Range("G19").Value = 1.103000044822693
$ws.Range("H19").Value = 26730988
$ws.Range("I19").Value = "MDAI"
$ws.Range("D20").Value = 2.539999961853028
$ws.Range("E20").Value = 2.480000019073486
$ws.Range("F20").Value = 3.210000038146973
$ws.Range("G20").Value = 2.349999904632568
$ws.Range("H20").Value = 26730988
$ws.Range("I20").Value = "MDAI"
